$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 10 data rows (rows 2-11) for weekly price reports.
# A new week of data is being added on top (rows 2-3), and all the existing
# rows shift down by two (old row 2 -> new row 4, ..., old row 11 -> new row 13).
#
# To avoid Excel's row-insert picking up the bold header formatting (row 1),
# we insert the two new blank rows at the bottom of the existing data (after
# row 11, which has plain/un-bolded formatting) and then copy the old data
# down into its new shifted location. Finally we overwrite the top two rows
# with the new week's values.

# Step 1: insert 2 blank rows after the last data row (row 11), inheriting the
# plain (non-bold) formatting/number formats already used by that row.
$ws.Rows("12:13").Insert()

# Step 2: shift the existing 10 data rows (2-11) down to rows 4-13.
$ws.Range("A2:T11").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Step 3: write the new week's data into rows 2-3. These reuse the same
# market/product/quality/unit metadata as the (former) first two rows,
# now located at rows 4-5, just with an updated date and prices.
$ws.Range("A2").Value2 = 11
$ws.Range("B2").Value2 = "Vega Monumental Concepción"
$ws.Range("C2").Value2 = "Bíobío"
$ws.Range("D2").Value2 = 44602
$ws.Range("E2").Value2 = 8
$ws.Range("F2").Value2 = "Fruta"
$ws.Range("G2").Value2 = 100101
$ws.Range("H2").Value2 = "Berries"
$ws.Range("I2").Value2 = 100101004
$ws.Range("J2").Value2 = "Frambuesa"
$ws.Range("K2").Value2 = "Sin especificar"
$ws.Range("L2").Value2 = "Primera"
$ws.Range("M2").Value2 = 200
$ws.Range("N2").Value2 = 6000
$ws.Range("O2").Value2 = 7000
$ws.Range("P2").Value2 = 6500
$ws.Range("Q2").Value2 = "$/bandeja 2 kilos"
$ws.Range("R2").Value2 = "Región de Ñuble"
$ws.Range("S2").Value2 = 3250
$ws.Range("T2").Value2 = 2

$ws.Range("A3").Value2 = 11
$ws.Range("B3").Value2 = "Vega Monumental Concepción"
$ws.Range("C3").Value2 = "Bíobío"
$ws.Range("D3").Value2 = 44602
$ws.Range("E3").Value2 = 8
$ws.Range("F3").Value2 = "Fruta"
$ws.Range("G3").Value2 = 100101
$ws.Range("H3").Value2 = "Berries"
$ws.Range("I3").Value2 = 100101004
$ws.Range("J3").Value2 = "Frambuesa"
$ws.Range("K3").Value2 = "Sin especificar"
$ws.Range("L3").Value2 = "Segunda"
$ws.Range("M3").Value2 = 100
$ws.Range("N3").Value2 = 5000
$ws.Range("O3").Value2 = 5000
$ws.Range("P3").Value2 = 5000
$ws.Range("Q3").Value2 = "$/bandeja 2 kilos"
$ws.Range("R3").Value2 = "Región de Ñuble"
$ws.Range("S3").Value2 = 2500
$ws.Range("T3").Value2 = 2

Write-Output "Applied weekly update: rows shifted, new week added at top."
